$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number format of the Fecha column (D) from the row above so the
# new date cell keeps the same "YYYY-MM-DD HH:MM:SS" style as the rest of
# the column.
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat

$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44511
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 300000000
$ws.Cells.Item(7, 7).Value = "Espárragos"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 600
$ws.Cells.Item(7, 11).Value = 1300
$ws.Cells.Item(7, 12).Value = 1400
$ws.Cells.Item(7, 13).Value = 1350
$ws.Cells.Item(7, 14).Value = "$/kilo"
$ws.Cells.Item(7, 15).Value = "Provincia de Linares"
$ws.Cells.Item(7, 16).Value = 1350
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
